$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 8, pushing existing rows 8-54 down to 9-55.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with this week's data. Columns A, B, C,
# E, F, G, H and R are constant for every row in this sheet, so reuse the
# values already present in (now) row 9, which used to be row 8.
$ws.Range("A8").Value = 11
$ws.Range("B8").Value = "Vega Monumental Concepción"
$ws.Range("C8").Value = "Bíobío"
$ws.Range("D8").Value = 44750
$ws.Range("D8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 100112037
$ws.Range("G8").Value = "Cebollín"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 130
$ws.Range("K8").Value = 7500
$ws.Range("L8").Value = 8000
$ws.Range("M8").Value = 7808
$ws.Range("N8").Value = "`$/paquete 36 unidades"
$ws.Range("O8").Value = "Región Metropolitana"
$ws.Range("P8").Value = 217
$ws.Range("Q8").Value = 36
$ws.Range("R8").Value = "Hortaliza"
